$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alascca table of mutation class")
$ws.Columns.Item(6).ColumnWidth = 0
$ws.Columns.Item(7).ColumnWidth = 1
$ws.Columns.Item(8).ColumnWidth = 100
Write-Host "done"
